# Modify the original data table style: re-order the 3 "Δθ/Δq n" / "q n"
# column pairs so each "q n" column sits directly to the right of its
# matching "Δθ/Δq n" column. Column A (Δθ/Δq 1) and F (q 3) stay put; the
# B/C/D/E columns are cyclically rotated:
#   new B = old D (q 1)
#   new C = old B (Δθ/Δq 2)
#   new D = old E (q 2)
#   new E = old C (Δθ/Δq 3)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells(1, 1).End(4).Row  # xlDown = 4 (from A1 down through the used column)
if ($lastRow -lt 10) { $lastRow = 10 }

for ($r = 1; $r -le $lastRow; $r++) {
    $bVal = $ws.Cells.Item($r, 2).Value2
    $cVal = $ws.Cells.Item($r, 3).Value2
    $dVal = $ws.Cells.Item($r, 4).Value2
    $eVal = $ws.Cells.Item($r, 5).Value2

    $ws.Cells.Item($r, 2).Value = $dVal
    $ws.Cells.Item($r, 3).Value = $bVal
    $ws.Cells.Item($r, 4).Value = $eVal
    $ws.Cells.Item($r, 5).Value = $cVal
}
